$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update site name and comment for the Karasjok row (row 14)
$ws.Range("A14").Value = "Karasjok Camping"
$ws.Range("E14").Value = "Suggested by Ann Kristin from Helitrans in e-mail received 02.09.2021"

# Update coordinates supplied by Helitrans
$ws.Range("B14").Value = 69.467986997867399
$ws.Range("C14").Value = 25.487036705017001

# Update the active selection cell as in the saved workbook
$ws.Range("C19").Select()
